# Update the "想去人数" (want-to-go count) figures that changed between
# crawler runs, on both the "展览" sheet and the aggregated "全部类型" sheet.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 388
$ws1.Range("F3").Value = 681
$ws1.Range("F7").Value = 10986
$ws1.Range("F12").Value = 10841
$ws1.Range("F15").Value = 25
$ws1.Range("F16").Value = 755
$ws1.Range("F17").Value = 5438
$ws1.Range("F18").Value = 84
$ws1.Range("F19").Value = 3406

# Sheet "全部类型" (all types, aggregated view)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 388
$ws4.Range("F3").Value = 681
$ws4.Range("F10").Value = 10986
$ws4.Range("F15").Value = 10841
$ws4.Range("F18").Value = 25
$ws4.Range("F19").Value = 755
$ws4.Range("F20").Value = 5438
$ws4.Range("F21").Value = 84
$ws4.Range("F22").Value = 3406
